$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 4500
$ws.Range("I74").Value = 4500
$ws.Range("K74").Value = 4500
$ws.Range("M74").Value = -3564
# Row 77
$ws.Range("H77").Value = 4500
$ws.Range("I77").Value = 4500
$ws.Range("K77").Value = 22500
$ws.Range("M77").Value = -17820
# Row 98
$ws.Range("H98").Value = 2571.75
$ws.Range("I98").Value = 2571.75
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2571.75
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -1073.75
# Row 113
$ws.Range("H113").Value = 2032.5
$ws.Range("I113").Value = 811.75
$ws.Range("J113").Value = 4474
$ws.Range("K113").Value = 811.75
$ws.Range("L113").Value = 4474
$ws.Range("M113").Value = 2442.25
$ws.Range("N113").Value = -10982
# Row 122
$ws.Range("H122").Value = 2571.75
$ws.Range("I122").Value = 2571.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7715.25
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -5265.25
# Row 129
$ws.Range("H129").Value = 1933
$ws.Range("I129").Value = 2119.6
$ws.Range("J129").Value = 1000
$ws.Range("K129").Value = 6358.799999999999
$ws.Range("L129").Value = 3000
$ws.Range("M129").Value = -1358.799999999999
$ws.Range("N129").Value = -13000
# Row 132
$ws.Range("H132").Value = 3522.75
$ws.Range("I132").Value = 3497.4443
$ws.Range("K132").Value = 10492.3329
$ws.Range("M132").Value = -7962.332900000001
# Row 135
$ws.Range("H135").Value = 6146.826
$ws.Range("I135").Value = 1745.1666
$ws.Range("J135").Value = 21992.8
$ws.Range("K135").Value = 15706.4994
$ws.Range("L135").Value = 197935.2
$ws.Range("M135").Value = -13171.4994
$ws.Range("N135").Value = -203005.2

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 2811544.2
$ws.Range("I6").Value = 15972.5
$ws.Range("J6").Value = 5048002
$ws.Range("K6").Value = 15972.5
$ws.Range("L6").Value = 5048002
$ws.Range("M6").Value = -15799.5
$ws.Range("N6").Value = -5048348
# Row 61
$ws.Range("H61").Value = 5712.8887
$ws.Range("I61").Value = 5632.8335
$ws.Range("K61").Value = 5632.8335
$ws.Range("M61").Value = -5420.8335
# Row 74
$ws.Range("H74").Value = 2961.8333
$ws.Range("I74").Value = 2822
$ws.Range("K74").Value = 2822
$ws.Range("M74").Value = -1948
# Row 77
$ws.Range("H77").Value = 2961.8333
$ws.Range("I77").Value = 2822
$ws.Range("K77").Value = 14110
$ws.Range("M77").Value = -9742
# Row 132
$ws.Range("H132").Value = 5844.6665
$ws.Range("I132").Value = 5775.8423
$ws.Range("K132").Value = 17327.5269
$ws.Range("M132").Value = -14797.5269
# Row 136
$ws.Range("H136").Value = 5712.8887
$ws.Range("I136").Value = 5632.8335
$ws.Range("K136").Value = 16898.5005
$ws.Range("M136").Value = -14348.5005
# Row 138
$ws.Range("H138").Value = 88407.336
$ws.Range("J138").Value = 88407.336
$ws.Range("L138").Value = 88407.336
$ws.Range("N138").Value = -98687.336

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 10002450
$ws.Range("J7").Value = 16670000
$ws.Range("L7").Value = 16670000
$ws.Range("N7").Value = -16670226
# Row 94
$ws.Range("H94").Value = 1144.7646
$ws.Range("I94").Value = 1100.7333
$ws.Range("K94").Value = 1100.7333
$ws.Range("M94").Value = -649.7333000000001
# Row 107
$ws.Range("H107").Value = 3210.0476
$ws.Range("I107").Value = 3101.2778
$ws.Range("K107").Value = 3101.2778
$ws.Range("M107").Value = -1181.2778
# Row 134
$ws.Range("H134").Value = 3551.348
$ws.Range("I134").Value = 3654.7896
$ws.Range("K134").Value = 10964.3688
$ws.Range("M134").Value = -8429.3688

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 1028.7142
$ws.Range("I6").Value = 1116.8334
$ws.Range("K6").Value = 1116.8334
$ws.Range("M6").Value = -1003.8334
# Row 12
$ws.Range("H12").Value = 1578.8
$ws.Range("I12").Value = 1348.5
$ws.Range("K12").Value = 1348.5
$ws.Range("M12").Value = -1178.5
# Row 31
$ws.Range("H31").Value = 5853.875
$ws.Range("I31").Value = 6529
$ws.Range("K31").Value = 6529
$ws.Range("M31").Value = -6234
# Row 34
$ws.Range("H34").Value = 5853.875
$ws.Range("I34").Value = 6529
$ws.Range("K34").Value = 6529
$ws.Range("M34").Value = -6327
# Row 62
$ws.Range("H62").Value = 2955.2
$ws.Range("I62").Value = 2929.3333
$ws.Range("J62").Value = 2994
$ws.Range("K62").Value = 2929.3333
$ws.Range("L62").Value = 2994
$ws.Range("M62").Value = -2305.3333
$ws.Range("N62").Value = -4242
# Row 65
$ws.Range("H65").Value = 2955.2
$ws.Range("I65").Value = 2929.3333
$ws.Range("J65").Value = 2994
$ws.Range("K65").Value = 14646.6665
$ws.Range("L65").Value = 14970
$ws.Range("M65").Value = -11526.6665
$ws.Range("N65").Value = -21210
# Row 132
$ws.Range("H132").Value = 6754.8423
$ws.Range("I132").Value = 5421.696
$ws.Range("J132").Value = 8799
$ws.Range("K132").Value = 16265.088
$ws.Range("L132").Value = 26397
$ws.Range("M132").Value = -13735.088
$ws.Range("N132").Value = -31457

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 56
$ws.Range("H56").Value = 7481.049
$ws.Range("I56").Value = 7481.049
$ws.Range("K56").Value = 7481.049
$ws.Range("M56").Value = -6951.049
# Row 113
$ws.Range("H113").Value = 1990.75
$ws.Range("J113").Value = 2499
$ws.Range("L113").Value = 7497
$ws.Range("N113").Value = -11837
# Row 132
$ws.Range("H132").Value = 1699.3334
$ws.Range("I132").Value = 799.25
$ws.Range("J132").Value = 3499.5
$ws.Range("K132").Value = 7193.25
$ws.Range("L132").Value = 31495.5
$ws.Range("M132").Value = -4663.25
$ws.Range("N132").Value = -36555.5
# Row 136
$ws.Range("H136").Value = 16572.2
$ws.Range("I136").Value = 17765.25
$ws.Range("K136").Value = 53295.75
$ws.Range("M136").Value = -48195.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 10
$ws.Range("H10").Value = 18001
$ws.Range("I10").Value = 4000
$ws.Range("J10").Value = 25001.5
$ws.Range("K10").Value = 4000
$ws.Range("L10").Value = 25001.5
$ws.Range("M10").Value = -3831
$ws.Range("N10").Value = -25339.5
# Row 47
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").ClearContents()
$ws.Range("N47").Value = 0
# Row 52
$ws.Range("H52").Value = 30136
$ws.Range("I52").Value = 27272
$ws.Range("K52").Value = 27272
$ws.Range("M52").Value = -27013
# Row 93
$ws.Range("H93").Value = 39999
$ws.Range("J93").Value = 39999
$ws.Range("L93").Value = 39999
$ws.Range("N93").Value = -43743
# Row 107
$ws.Range("H107").Value = 647.5
$ws.Range("I107").Value = 597.1429000000001
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 597.1429000000001
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1322.8571
$ws.Range("N107").Value = -4840
# Row 132
$ws.Range("H132").Value = 7290.9546
$ws.Range("I132").Value = 6994
$ws.Range("J132").Value = 8082.8335
$ws.Range("K132").Value = 20982
$ws.Range("L132").Value = 24248.5005
$ws.Range("M132").Value = -18452
$ws.Range("N132").Value = -29308.5005

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 3604
$ws.Range("I22").Value = 3388.4
$ws.Range("J22").Value = 3873.5
$ws.Range("K22").Value = 3388.4
$ws.Range("L22").Value = 3873.5
$ws.Range("M22").Value = -3093.4
$ws.Range("N22").Value = -4463.5
# Row 27
$ws.Range("H27").Value = 3604
$ws.Range("I27").Value = 3388.4
$ws.Range("J27").Value = 3873.5
$ws.Range("K27").Value = 3388.4
$ws.Range("L27").Value = 3873.5
$ws.Range("M27").Value = -3281.4
$ws.Range("N27").Value = -4087.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 2492.5
$ws.Range("I107").Value = 2272.9048
$ws.Range("K107").Value = 6818.714399999999
$ws.Range("M107").Value = -4898.714399999999
# Row 136
$ws.Range("H136").Value = 2692.3147
$ws.Range("I136").Value = 2425.7
$ws.Range("K136").Value = 7277.099999999999
$ws.Range("M136").Value = -4727.099999999999
